$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the value for D2 (Mon, Wed column -> new energy value)
$ws.Range("D2").Value = 9

# Update the selected cell/range to D3
$ws.Range("D3").Select()
